$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data (header stays the same; rows 2-18 hold the player table)
$data = @(
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("Payton Pritchard",   "PG",       "Boston Celtics"),
    @("Bradley Beal",       "PG,SG,SF", "Phoenix Suns"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("De'Andre Hunter",    "SF,PF",    "Atlanta Hawks"),
    @("Victor Wembanyama",  "C",        "San Antonio Spurs"),
    @("Goga Bitadze",       "C",        "Orlando Magic"),
    @("Deandre Ayton",      "C",        "Portland Trail Blazers"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("Myles Turner",       "C",        "Indiana Pacers"),
    @("Jamal Murray",       "PG,SG",    "Denver Nuggets"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
    @("Tari Eason",         "SF,PF",    "Houston Rockets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
